# Atualização de bases das ligas, do dia: 15-06-2024 às 21:10
#
# The "id" (match id) / odds rows for a few fixtures were attached to the
# wrong match metadata. This swaps the full data payload (columns B..AD:
# id, Div, Date, HomeTeam, AwayTeam, scores, odds, etc.) between the
# affected rows while leaving column A (the sequential row index) fixed in
# place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 175 <-> 176 -------------------------------------------------
$row175 = $ws.Range("B175:AD175").Value2
$row176 = $ws.Range("B176:AD176").Value2
$ws.Range("B175:AD175").Value2 = $row176
$ws.Range("B176:AD176").Value2 = $row175

# --- Rows 183 <-> 184 -------------------------------------------------
$row183 = $ws.Range("B183:AD183").Value2
$row184 = $ws.Range("B184:AD184").Value2
$ws.Range("B183:AD183").Value2 = $row184
$ws.Range("B184:AD184").Value2 = $row183

# --- Rows 185 <-> 186 -------------------------------------------------
$row185 = $ws.Range("B185:AD185").Value2
$row186 = $ws.Range("B186:AD186").Value2
$ws.Range("B185:AD185").Value2 = $row186
$ws.Range("B186:AD186").Value2 = $row185

# --- Rows 313 -> 314 -> 315 -> 313 (3-way rotation) --------------------
$row313 = $ws.Range("B313:AD313").Value2
$row314 = $ws.Range("B314:AD314").Value2
$row315 = $ws.Range("B315:AD315").Value2
$ws.Range("B313:AD313").Value2 = $row314
$ws.Range("B314:AD314").Value2 = $row315
$ws.Range("B315:AD315").Value2 = $row313
